$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# First, relocate the existing "Bold outputs" note from C28 down to C31,
# freeing up row 28 (and 27) for the two new test rows being inserted.
$noteValue = $ws.Range("C28").Value()
$ws.Range("C31").Value = $noteValue
$ws.Range("C31").Font.Bold = $true

$ws.Range("C28").ClearContents()
$ws.Range("C28").Font.Bold = $false

# New rows 27-28 (values entered in this order so that the workbook's
# shared-string table is built up in the same sequence as the source edit)
$ws.Range("A27").Value = "player1 joins a game"
$ws.Range("C27").Value = "player1 sees the list of online players"
$ws.Range("A28").Value = "player24 logs into the server and enters ""player24"
$ws.Range("B27").Value = "player1 joins a game and still sees the list of online players"
$ws.Range("B28").Value = "player1 should see the newly added player"
$ws.Range("C28").Value = "player1 cannot see player24 on the online player list"
$ws.Range("C28").Font.Bold = $true

# Restore view state: scroll so row 8 is the topmost visible row, and
# select C22 as the active cell.
$ws.Range("C22").Select()
$excel.ActiveWindow.ScrollRow = 8
